$d = $word.ActiveDocument

function Replace-ParagraphWithXml {
    # Finds searchText (which must span the ENTIRE paragraph content, start to end,
    # because InsertXML with a <w:p> wrapper misbehaves -- it always appends new
    # content at the end of the target paragraph rather than at the deletion point
    # -- unless the replaced range is exactly the whole paragraph) and replaces the
    # paragraph's run content with innerXml.
    param(
        [string]$searchText,
        [string]$innerXml
    )
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $searchText"
        return
    }
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# --- 1. Heading: "Additional documentation for PavementEye Project" ---
Replace-ParagraphWithXml "Additional documentation for PavementEye Project" (
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Additional documentation for </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/><w:b/><w:bCs/></w:rPr><w:t>PavementEye</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Project</w:t></w:r>'
)

# --- 2. KPI table: "The system can achieve accuracy of mAP of 0.5 or higher in Egyptian roads " (whole paragraph) ---
Replace-ParagraphWithXml "The system can achieve accuracy of mAP of 0.5 or higher in Egyptian roads " (
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t xml:space="preserve">The system can achieve accuracy of </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t>mAP</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t xml:space="preserve"> of 0.5 or higher in Egyptian roads </w:t></w:r>'
)

# --- 3. KPI table: "Achieved" + " with mAP of 0.55..." (whole paragraph, incl. bold "Achieved" run) ---
Replace-ParagraphWithXml "Achieved with mAP of 0.55 after fine tuning on EGY_PDD dataset with more than 15K images of Egyptian roads." (
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/><w:b/><w:bCs/></w:rPr><w:t>Achieved</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t xml:space="preserve"> with </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t>mAP</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t xml:space="preserve"> of 0.55 after fine tuning on EGY_PDD dataset with more than 15K images of Egyptian roads.</w:t></w:r>'
)

# --- 4. KPI table: "The system can handle 1 image per second..." -> "...every 5 seconds..." (whole paragraph) ---
Replace-ParagraphWithXml "The system can handle 1 image per second for the whole pipeline cycle. From capturing till image storage." (
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t xml:space="preserve">The system can handle 1 image </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t xml:space="preserve">every 5 seconds </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t>for the whole pipeline cycle. From capturing till image storage.</w:t></w:r>'
)

# --- 5. KPI table: "Achieved" + " successfully and validated via Kafka UI..." (whole paragraph) ---
Replace-ParagraphWithXml "Achieved successfully and validated via Kafka UI and spark logging. Also enhanced more by replacing normal APIs with web sockets for duplex connection between backend and flutter application." (
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/><w:b/><w:bCs/></w:rPr><w:t>Achieved</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t xml:space="preserve"> successfully and validated </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t xml:space="preserve">(less than 5 seconds per image, most are between 2 and 3 seconds) </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t>via Kafka UI and spark logging. Also enhanced more by replacing normal APIs with web sockets for duplex connection between backend and flutter application.</w:t></w:r>'
)

# --- 6. KPI table: "... district (منتزة, رمل, عجمي, الخ)" -> split "منتزة" out with spellcheck markers (whole paragraph) ---
Replace-ParagraphWithXml "The system can correctly get location (latitude and longitude) and assign the crack accurately to its roads and district (منتزة, رمل, عجمي, الخ)" (
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t>The system can correctly get location (latitude and longitude) and assign the crack accurately to its roads and district (</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/><w:rtl/><w:lang w:bidi="ar-EG"/></w:rPr><w:t>منتزة</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/><w:rtl/><w:lang w:bidi="ar-EG"/></w:rPr><w:t>, رمل, عجمي, الخ</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t>)</w:t></w:r>'
)

# --- 7. "Achieved" + ". We added part on streamlit dashboard..." (whole paragraph) ---
Replace-ParagraphWithXml "Achieved. We added part on streamlit dashboard to type the name of the image and get it from the data lake and put the bounding boxes on it." (
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/><w:b/><w:bCs/></w:rPr><w:t>Achieved</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t xml:space="preserve">. We added part on </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t>streamlit</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Google Sans" w:hAnsi="Google Sans"/></w:rPr><w:t xml:space="preserve"> dashboard to type the name of the image and get it from the data lake and put the bounding boxes on it.</w:t></w:r>'
)

# --- 8. "Team Leader, Data Engineering (Spark, Cassandra, kafka), Docker, ... Streamlit dashboard." (whole paragraph) ---
Replace-ParagraphWithXml "Team Leader, Data Engineering (Spark, Cassandra, kafka), Docker, Computer Vision (fine tuning on EGY_PDD dataset), Cloud Computing (Azure VM and Datalake and Huawei Cloud ECS and OBS), Backend Development and websocket connection, Flutter Development, Streamlit dashboard." (
    '<w:r><w:t xml:space="preserve">Team Leader, Data Engineering (Spark, Cassandra, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>kafka</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>),</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Docker,</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Computer Vision (fine tuning on EGY_PDD dataset), Cloud Computing (Azure VM and </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Datalake</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> and Huawei Cloud ECS and OBS), Backend Development and </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>websocket</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> connection, Flutter Development, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Streamlit</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> dashboard</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>'
)

# --- 9. "Salsabel" -> wrap with spellcheck markers (whole paragraph) ---
Replace-ParagraphWithXml "Salsabel" (
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Salsabel</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)

# --- 10. "Huawei Cloud, Data visualization and streamlit dashboard, helping in testcases, Video editing." (whole paragraph) ---
Replace-ParagraphWithXml "Huawei Cloud, Data visualization and streamlit dashboard, helping in testcases, Video editing." (
    '<w:r><w:t xml:space="preserve">Huawei Cloud, Data visualization and </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>streamlit</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> dashboard, helping in testcases, Video editing</w:t></w:r>' +
    '<w:r><w:t>, AI Chatbot</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>'
)

Write-Output "All replacements done"
